$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.918.78'
$ws.Range("E2").Value = '  -0.34%  '

$ws.Range("D3").Value = '2.035.28'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'244.88"
$ws.Range("E5").Value = '  -1.35%  '

$ws.Range("D6").Value = "'0.657"

$ws.Range("D7").Value = "'57.80"
$ws.Range("E7").Value = '  +0.40%  '

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = "'0.375"
$ws.Range("E9").Value = '  -1.02%  '

$ws.Range("D10").Value = "'0.0766"
$ws.Range("E10").Value = '  -1.21%  '

$ws.Range("E11").Value = '  +2.31%  '

$ws.Range("D12").Value = "'15.34"
$ws.Range("E12").Value = '  -2.27%  '

$ws.Range("D13").Value = "'0.879"
$ws.Range("E13").Value = '  +8.94%  '

$ws.Range("D14").Value = '2.331.15'
$ws.Range("E14").Value = '  -0.62%  '

$ws.Range("D15").Value = "'5.62"
$ws.Range("E15").Value = '  +2.09%  '

$ws.Range("D16").Value = '2.009.38'
$ws.Range("E16").Value = '  -1.80%  '

$ws.Range("D17").Value = "'18.11"
$ws.Range("E17").Value = '  +8.93%  '

$ws.Range("D18").Value = '36.863.57'
$ws.Range("E18").Value = '  -0.46%  '

$ws.Range("D19").Value = "'73.35"
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").Value = '0.0₃0884'
$ws.Range("E20").Value = '  -1.21%  '

$ws.Range("D21").Value = "'5.35"
$ws.Range("E21").Value = '  +0.75%  '

$ws.Range("D22").Value = "'235.31"
$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("E24").Value = '  +3.38%  '

$ws.Range("D25").Value = "'9.58"
$ws.Range("E25").Value = '  +4.96%  '

$ws.Range("D26").Value = "'168.93"
$ws.Range("E26").Value = '  +0.84%  '

$ws.Range("E27").Value = '  -5.05%  '

$ws.Range("D28").Value = "'19.87"
$ws.Range("E28").Value = '  +0.91%  '

$ws.Range("D29").Value = "'5.44"
$ws.Range("E29").Value = '  +16.59%  '

$ws.Range("E30").Value = '  -0.37%  '

$ws.Range("D31").Value = "'1.10"
$ws.Range("E31").Value = '  -2.00%  '

$ws.Range("D32").Value = "'4.71"
$ws.Range("E32").Value = '  +6.63%  '

$ws.Range("D33").Value = "'0.0610"
$ws.Range("E33").Value = '  +0.07%  '

$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").Value = "'0.0864"
$ws.Range("E35").Value = '  -3.96%  '

$ws.Range("E36").Value = '  +6.34%  '

$ws.Range("D37").Value = "'2.23"
$ws.Range("E37").Value = '  +0.93%  '

$ws.Range("D38").Value = "'1.30"
$ws.Range("E38").Value = '  -3.47%  '

$ws.Range("E39").Value = '  -1.42%  '

$ws.Range("D40").Value = "'5.11"
$ws.Range("E40").Value = '  +0.95%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = "'0.0221"
$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").Value = "'0.0966"
$ws.Range("E42").Value = '  -10.42%  '

$ws.Range("E43").Value = '  +0.84%  '

$ws.Range("D44").Value = "'96.88"
$ws.Range("E44").Value = '  +1.77%  '

$ws.Range("D45").Value = "'16.85"
$ws.Range("E45").Value = '  -2.22%  '

$ws.Range("D46").Value = '1.291.04'
$ws.Range("E46").Value = '  +1.23%  '

$ws.Range("D47").Value = "'3.80"
$ws.Range("E47").Value = '  +10.93%  '

$ws.Range("D48").Value = "'2.33"
$ws.Range("E48").Value = '  -3.67%  '

$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = '  -0.37%  '

$ws.Range("D50").Value = "'6.68"
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").Value = '2.218.98'
$ws.Range("E51").Value = '  -0.64%  '
